# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Ají".
# The new record is inserted at row 259, pushing the former rows 259-275 down to 260-276.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(259).Insert()

$ws.Cells.Item(259, 1).Value = 4
$ws.Cells.Item(259, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(259, 3).Value = "Los Lagos"
$ws.Cells.Item(259, 4).Value = 44746
$ws.Cells.Item(259, 5).Value = 10
$ws.Cells.Item(259, 6).Value = 100112021
$ws.Cells.Item(259, 7).Value = "Ají"
$ws.Cells.Item(259, 8).Value = "Inferno"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 50
$ws.Cells.Item(259, 11).Value = 19000
$ws.Cells.Item(259, 12).Value = 19000
$ws.Cells.Item(259, 13).Value = 19000
$ws.Cells.Item(259, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(259, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(259, 16).Value = 1583
$ws.Cells.Item(259, 17).Value = 12
$ws.Cells.Item(259, 18).Value = "Hortaliza"
